$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45243 = 2023-11-13) for
# every data row (2-18). The automatic update bumps each of these by one day
# (serial 45244 = 2023-11-14).
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
